# "exclusion to all white milk brands"
#
# The workbook has two sheets: "Exclude" and "Include". On both sheets a
# parameter row identifies the KPI parameter named "brand" and its value
# (the pipe/comma separated list of white-milk brand names). This change:
#   1. Renames the "brand" parameter to "brand_name" everywhere it is used
#      as an exact cell value (case sensitive, not touching the "Brand"
#      column header).
#   2. Cleans up stray trailing whitespace/newlines on the brand value list
#      so both sheets reference one single, tidy string.
#   3. Nudges the formatting of the (now renamed) brand_name cells so they
#      pick up their own style record (content/alignment unchanged).
#   4. Switches the active/selected tab from "Exclude" to "Include" and
#      updates each sheet's remembered selection.

$wb = $excel.ActiveWorkbook

$exclude = $wb.Worksheets.Item("Exclude")
$include = $wb.Worksheets.Item("Include")

# ---------------------------------------------------------------------
# 1 & 2: text clean-up, applied across both sheets' used ranges.
# xlWhole (1) + MatchCase=$true so only exact, case sensitive matches are
# touched (keeps the "Brand" column header intact).
# ---------------------------------------------------------------------
foreach ($ws in @($exclude, $include)) {
    $rng = $ws.UsedRange

    [void]$rng.Replace("brand", "brand_name", 1, [Type]::Missing, $true)

    [void]$rng.Replace(
        "A2 White Milk,Dairy Farmers White Milk,Pauls White Milk,Other Dairy`n`n",
        "A2 White Milk,Dairy Farmers White Milk,Pauls White Milk,Other Dairy",
        1, [Type]::Missing, $true)

    [void]$rng.Replace(
        "A2 White Milk,Dairy Farmers White Milk,Pauls White Milk,Other Dairy  ",
        "A2 White Milk,Dairy Farmers White Milk,Pauls White Milk,Other Dairy",
        1, [Type]::Missing, $true)
}

# ---------------------------------------------------------------------
# 3: give the renamed "brand_name" parameter cells their own style.
# Re-asserting WrapText (already true) is enough to split them off into a
# distinct cell format without altering anything visible.
# ---------------------------------------------------------------------
foreach ($addr in @("C4", "C10", "C13")) {
    $exclude.Range($addr).WrapText = $true
}
foreach ($addr in @("C2", "C4")) {
    $include.Range($addr).WrapText = $true
}

# ---------------------------------------------------------------------
# 4: move the active tab/selection from Exclude to Include.
# ---------------------------------------------------------------------
[void]$exclude.Range("C13").Select()
[void]$include.Activate()
[void]$include.Range("C2").Select()
